$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.723.61"
$ws.Cells.Item(2, 5).Value = "  +0.34%  "
$ws.Cells.Item(3, 4).Value = "1.601.71"
$ws.Cells.Item(3, 5).Value = "  +0.31%  "
$ws.Cells.Item(4, 5).Value = "  +0.21%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "211.34"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.04%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.513"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.54%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.01"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.23%  "
$ws.Cells.Item(8, 5).Value = "  +0.21%  "
$ws.Cells.Item(9, 5).Value = "  +0.58%  "
$ws.Cells.Item(10, 5).Value = "  +1.19%  "
$ws.Cells.Item(11, 5).Value = "  +0.79%  "
$ws.Cells.Item(12, 4).Value = "1.826.21"
$ws.Cells.Item(13, 4).Value = "1.601.75"
$ws.Cells.Item(13, 5).Value = "  -0.29%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.04"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.43%  "
$ws.Cells.Item(15, 5).Value = "  +0.45%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "65.20"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.50%  "
$ws.Cells.Item(17, 4).Value = "26.697.05"
$ws.Cells.Item(17, 5).Value = "  +0.28%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0746"
$ws.Cells.Item(18, 5).Value = "  +1.28%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "210.82"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.13%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.21"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.46%  "
$ws.Cells.Item(21, 5).Value = "  +0.14%  "
$ws.Cells.Item(22, 5).Value = "  +0.85%  "
$ws.Cells.Item(23, 5).Value = "  +0.42%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.99"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.05%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "143.62"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.14%  "
$ws.Cells.Item(26, 5).Value = "  +0.36%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.12"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.38%  "
$ws.Cells.Item(28, 5).Value = "  -0.75%  "
$ws.Cells.Item(29, 5).Value = "  +1.10%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0514"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.34%  "
$ws.Cells.Item(31, 5).Value = "  -0.29%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.26"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.33%  "
$ws.Cells.Item(33, 5).Value = "  +1.67%  "
$ws.Cells.Item(34, 4).Value = "1.298.28"
$ws.Cells.Item(34, 5).Value = "  +1.77%  "
$ws.Cells.Item(35, 5).Value = "  +0.56%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.610"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.52%  "
$ws.Cells.Item(37, 5).Value = "  +0.82%  "
$ws.Cells.Item(38, 5).Value = "  +20.47%  "
$ws.Cells.Item(40, 5).Value = "  -1.94%  "
$ws.Cells.Item(41, 5).Value = "  -1.41%  "
$ws.Cells.Item(42, 5).Value = "  -0.40%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.782"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.21%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "63.27"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.22%  "
$ws.Cells.Item(45, 4).Value = "1.737.09"
$ws.Cells.Item(45, 5).Value = "  +0.25%  "
$ws.Cells.Item(46, 5).Value = "  +0.95%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.57"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.01%  "
$ws.Cells.Item(48, 4).Value = "0.0₆0105"
$ws.Cells.Item(48, 5).Value = "  -1.28%  "
$ws.Cells.Item(49, 5).Value = "  -0.45%  "
$ws.Cells.Item(50, 5).Value = "  +1.90%  "
